$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "scroll into view"
$ws.Range("B7").Value = "<window_name>|<control_name(list item)>"

$ws.Range("A7").Font.Bold = $true

$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:B7"))

$ws.Range("B7").Select()
